$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.311.60'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '2.638.92'
$ws.Range('E3').Value = '  -3.36%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.33'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.69'
$ws.Range('E6').Value = '  -1.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.544'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('D9').Value = '2.639.40'
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.144'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.363'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.22'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.99'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').Value = '3.124.25'
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('D17').Value = '67.266.44'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').Value = '2.619.18'
$ws.Range('E18').Value = '  -4.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.89'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.87'
$ws.Range('E20').Value = '  +2.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '363.61'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.77'
$ws.Range('E23').Value = '  -3.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.08'
$ws.Range('E24').Value = '  +11.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  -5.55%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.98'
$ws.Range('E27').Value = '  -3.57%  '
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').Value = '  -3.36%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '554.22'
$ws.Range('E31').Value = '  -5.83%  '
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.38'
$ws.Range('E33').Value = '  -4.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.92'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.97'
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.41'
$ws.Range('E39').Value = '  -2.80%  '
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.26'
$ws.Range('E41').Value = '  -4.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.82'
$ws.Range('E42').Value = '  -5.20%  '
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.52'
$ws.Range('E45').Value = '  -5.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.10'
$ws.Range('E46').Value = '  -2.37%  '
$ws.Range('D47').Value = '0.0₆0302'
$ws.Range('E47').Value = '  -2.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.594'
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '153.90'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.88'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.73'
$ws.Range('E51').Value = '  -3.27%  '
